$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Capture the original column widths of D and E before anything else
# changes, so they can be swapped further down
$widthD = $ws.Columns.Item(4).ColumnWidth
$widthE = $ws.Columns.Item(5).ColumnWidth

# Waist size values that will occupy the new column D
$waist = @(36,25,28,35,36,30,27,37,40,32,48,42,36,38)

# Style values that currently live in column D need to move to column E
$styles = @()
for ($r = 2; $r -le 15; $r++) {
    $styles += $ws.Cells.Item($r, 4).Value()
}

# Update header row: D1 = "Waist Size", E1 = "Baggy Pants" (former D1 header)
$ws.Cells.Item(1, 5).Value = $ws.Cells.Item(1, 4).Value()
$ws.Cells.Item(1, 4).Value = "Waist Size"

# Move the old column D (style) values into column E, then put the waist
# size numbers into column D
for ($i = 0; $i -lt 14; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value = $styles[$i]
    $ws.Cells.Item($r, 4).Value = $waist[$i]
}

# Swap the custom column widths of D and E to match the moved content
$ws.Columns.Item(4).ColumnWidth = $widthE
$ws.Columns.Item(5).ColumnWidth = $widthD

# Update the selected cell to reflect the author's final cursor position
$ws.Range("F5").Select()
